$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the value for F15 (hours worked), which recalculates G15/H15 formulas.
$ws.Range("F15").Value = 20.25

# Update the active cell selection to F15
$ws.Range("F15").Select()
